$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.466.77"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "3.426.27"
$ws.Range("E3").Value = "  -2.64%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.90%  "
$ws.Range("D7").Value = "3.425.24"
$ws.Range("E7").Value = "  -2.64%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("E10").Value = "  -5.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.121"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -9.37%  "
$ws.Range("E12").Value = "  -7.27%  "
$ws.Range("D13").Value = "4.008.63"
$ws.Range("E13").Value = "  -2.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000180"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -10.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -8.36%  "
$ws.Range("D16").Value = "3.435.33"
$ws.Range("E16").Value = "  -1.96%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.115"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "65.431.81"
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("E19").Value = "  -10.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "391.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.57%  "
$ws.Range("E23").Value = "  -6.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.63%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "3.563.75"
$ws.Range("E26").Value = "  -2.65%  "
$ws.Range("E27").Value = "  -8.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -9.18%  "
$ws.Range("E31").Value = "  -8.99%  "
$ws.Range("D32").Value = "3.431.89"
$ws.Range("E32").Value = "  -2.43%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("E34").Value = "  -6.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "22.99"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "172.87"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.83"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -9.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.15"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -9.51%  "
$ws.Range("E39").Value = "  -7.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.53%  "
$ws.Range("E41").Value = "  -6.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.819"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "43.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -12.55%  "
$ws.Range("E46").Value = "  -9.68%  "
$ws.Range("E47").Value = "  +1.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -13.93%  "
$ws.Range("D51").Value = "2.192.94"
$ws.Range("E51").Value = "  -6.98%  "
